# Generate Report for Handoff
# Updates the "8fee2a14-c989-43aa-b604-5963f8c03a99.md" row's handoff
# timestamps across the Overview / zh-cn / de-de sheets, reflecting a
# freshly generated handoff xliff for that file.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: column G = "Latest HO Xliff Generate Date", row 7 is the
# 8fee2a14-c989-43aa-b604-5963f8c03a99.md file.
$overview.Range("G7").Value = "2016-08-26 10:40:16"

# zh-cn sheet: column H = "Latest Handoff Datetime", row 7 is the same file.
$zhcn.Range("H7").Value = "2016-08-26 10:40:00"

# de-de sheet: column H = "Latest Handoff Datetime", row 7 is the same file.
$dede.Range("H7").Value = "2016-08-26 10:40:16"
